$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.307.54'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.887.16'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = "'314.89"
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').Value = "'0.5147"
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').Value = "'0.3922"
$ws.Range('E8').Value = '  +3.32%  '
$ws.Range('D9').Value = "'0.08385"
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = "'1.125"
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('D11').Value = "'41.67"
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').Value = "'6.245"
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = "'20.78"
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.869.76'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').Value = "'7.306"
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = "'1.009"
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = "'91.69"
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').Value = "'0.06676"
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').Value = "'17.85"
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = "'6.062"
$ws.Range('E22').Value = '  +0.83%  '
$ws.Range('D23').Value = '28.339.96'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').Value = "'11.20"
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('D25').Value = "'2.277"
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = "'3.402"
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '2.098.76'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'2.531"
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = "'159.36"
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'20.72"
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = "'126.68"
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.1070"
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'1.051"
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'5.918"
$ws.Range('E34').Value = '  +6.21%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'3.611"
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = "'9.793"
$ws.Range('E36').Value = '  +2.61%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.02470"
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.06592"
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = "'0.2202"
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = "'1.217"
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.6559"
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').Value = "'5.034"
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'1.235"
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'11.35"
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.6177"
$ws.Range('E45').Value = '  +2.39%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'13.15"
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'1.290"
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = "'3.691"
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'2.023"
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = "'1.241"
$ws.Range('E50').Value = '  +3.02%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = "'121.75"
$ws.Range('E51').Value = '  +1.27%  '
